$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the B8:B78 range to text format so purely-numeric replacement
# strings are stored as text (matching the original inline-string cells)
# rather than being auto-coerced to numbers.
$ws.Range("B8:B78").NumberFormat = "@"

$ws.Range("B8").Value = "496"
$ws.Range("B9").Value = "324"
$ws.Range("B10").Value = "499"
$ws.Range("B11").Value = "502"
$ws.Range("B12").Value = "505"
$ws.Range("B13").Value = "508"
$ws.Range("B14").Value = "335"
$ws.Range("B15").Value = "511"
$ws.Range("B16").Value = "514"
$ws.Range("B17").Value = "517"
$ws.Range("B18").Value = "345"
$ws.Range("B19").Value = "520"
$ws.Range("B20").Value = "523"
$ws.Range("B21").Value = "526"
$ws.Range("B22").Value = "355"
$ws.Range("B23").Value = "529"
$ws.Range("B24").Value = "532"
$ws.Range("B25").Value = "535"
$ws.Range("B26").Value = "538"
$ws.Range("B27").Value = "367"
$ws.Range("B28").Value = "541"
$ws.Range("B29").Value = "544"
$ws.Range("B30").Value = "547"
$ws.Range("B31").Value = "377"
$ws.Range("B32").Value = "550"
$ws.Range("B33").Value = "553"
$ws.Range("B34").Value = "556"
$ws.Range("B35").Value = "387"
$ws.Range("B36").Value = "559"
$ws.Range("B37").Value = "562"
$ws.Range("B38").Value = "565"
$ws.Range("B39").Value = "397"
$ws.Range("B40").Value = "568"
$ws.Range("B41").Value = "571"
$ws.Range("B42").Value = "574"
$ws.Range("B43").Value = "577"
$ws.Range("B44").Value = "409"
$ws.Range("B45").Value = "580"
$ws.Range("B46").Value = "583"
$ws.Range("B47").Value = "586"
$ws.Range("B48").Value = "419"
$ws.Range("B49").Value = "589"
$ws.Range("B50").Value = "592"
$ws.Range("B51").Value = "595"
$ws.Range("B52").Value = "429"
$ws.Range("B53").Value = "598"
$ws.Range("B54").Value = "601"
$ws.Range("B55").Value = "604"
$ws.Range("B56").Value = "439"
$ws.Range("B57").Value = "607"
$ws.Range("B58").Value = "610"
$ws.Range("B59").Value = "613"
$ws.Range("B60").Value = "449"
$ws.Range("B61").Value = "616"
$ws.Range("B62").Value = "619"
$ws.Range("B63").Value = "622"
$ws.Range("B64").Value = "459"
$ws.Range("B65").Value = "625"
$ws.Range("B66").Value = "628"
$ws.Range("B67").Value = "631"
$ws.Range("B68").Value = "634"
$ws.Range("B69").Value = "471"
$ws.Range("B70").Value = "637"
$ws.Range("B71").Value = "640"
$ws.Range("B72").Value = "643"
$ws.Range("B73").Value = "481"
$ws.Range("B74").Value = "646"
$ws.Range("B75").Value = "649"
$ws.Range("B76").Value = "652"
$ws.Range("B77").Value = "491"
$ws.Range("B78").Value = "493"

# Restore default (unstyled) formatting so only the cell values change.
$ws.Range("B8:B78").Style = "Normal"

